# Results from R script
# - Correct the timestamp recorded for the 2024-07-08 bar (A100): the raw
#   pull had an intraday timestamp, normalize it to the 07:00 snapshot used
#   by every other row.
# - Append the newly scraped bar for 2024-07-09 (row 101).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix A100 ---------------------------------------------------------
$ws.Range("A100").Value = 45481.2916666667

# --- Append new row 101 ------------------------------------------------
# Column A (date): clone the date/time number format from the row above so
# the new cell keeps the same "yyyy-mm-dd hh:mm:ss" style, then set value.
$ws.Range("A100").Copy()
$ws.Range("A101").PasteSpecial(-4122)
$ws.Range("A101").Value = 45482.6292013889

$ws.Range("B101").Value = 6300
$ws.Range("C101").Value = 6.17999982833862
$ws.Range("D101").Value = 6.03999996185303
$ws.Range("E101").Value = 6.1399998664856
$ws.Range("F101").Value = 6.03999996185303

# Column G (adj_close) is stored as text in this sheet (matches every other
# row) - force text entry, then drop back to the default cell style so we
# don't leave a stray number format on the cell.
$ws.Range("G101").NumberFormat = "@"
$ws.Range("G101").Value = "6.03999996185303"
$ws.Range("G101").Style = "Normal"

# Column H (ticker)
$ws.Range("H101").Value = "PAL.MI"
